$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text before writing, so numeric-looking
# strings like "582.13" or "0.999" stay text (matching the source inlineStr cells)
# instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$originalStyle = $dataRange.Style
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.102.95"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "2.556.28"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "582.13"
$ws.Range("E5").Value = "  +1.72%  "

$ws.Range("D6").Value = "147.14"
$ws.Range("E6").Value = "  -2.24%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "5.58"
$ws.Range("E10").Value = "  -1.87%  "

$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").Value = "27.55"
$ws.Range("E13").Value = "  -1.93%  "

$ws.Range("D14").Value = "3.011.02"

$ws.Range("D15").Value = "62.993.97"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").Value = "2.552.37"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").Value = "11.37"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("D19").Value = "340.41"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("D20").Value = "4.37"
$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").Value = "6.78"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "65.85"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("D24").Value = "2.675.72"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "0.170"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "1.63"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  -3.48%  "

$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").Value = "8.37"
$ws.Range("E29").Value = "  -1.14%  "

$ws.Range("D30").Value = "7.82"
$ws.Range("E30").Value = "  +6.93%  "

$ws.Range("D31").Value = "1.97"
$ws.Range("E31").Value = "  +4.93%  "

$ws.Range("D32").Value = "0.0₃0820"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").Value = "177.32"
$ws.Range("E33").Value = "  -0.36%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "427.93"
$ws.Range("E35").Value = "  +0.66%  "

$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "0.405"
$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "19.18"
$ws.Range("E37").Value = "  +1.01%  "

$ws.Range("D38").Value = "4.41"
$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  -2.00%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "39.67"
$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("D43").Value = "151.11"
$ws.Range("E43").Value = "  -1.58%  "

$ws.Range("D44").Value = "3.79"
$ws.Range("E44").Value = "  +0.29%  "

$ws.Range("D45").Value = "20.91"
$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("D46").Value = "0.0545"
$ws.Range("E46").Value = "  +3.85%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.604"
$ws.Range("E47").Value = "  -1.23%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.0973"
$ws.Range("E48").Value = "  +0.65%  "

$ws.Range("E49").Value = "  +0.21%  "

$ws.Range("D50").Value = "18.29"
$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "1.72"
$ws.Range("E51").Value = "  -4.62%  "

# Restore the original (default) style so no stray number-format / style
# artifacts are left behind on the cells.
$dataRange.Style = $originalStyle

Write-Host "Updated cryptos list"